$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (shifts existing rows 26-71 down to 27-72)
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the new market-report entry
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44477
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100102
$ws.Range("H26").Value = "Cítricos"
$ws.Range("I26").Value = 100102004
$ws.Range("J26").Value = "Mandarina"
$ws.Range("K26").Value = "Murcott"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 270
$ws.Range("N26").Value = 11000
$ws.Range("O26").Value = 12000
$ws.Range("P26").Value = 11500
$ws.Range("Q26").Value = "$/caja 20 kilos"
$ws.Range("R26").Value = "Provincia de Melipilla"
$ws.Range("S26").Value = 575
$ws.Range("T26").Value = 20
